$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.939.71"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").Value = "3.569.47"
$ws.Range("E3").Value = "  -1.87%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'577.54"
$ws.Range("E5").Value = "  -2.88%  "
$ws.Range("D6").Value = "'188.65"
$ws.Range("E6").Value = "  -1.71%  "
$ws.Range("D7").Value = "'0.632"
$ws.Range("E7").Value = "  -2.65%  "
$ws.Range("D8").Value = "3.563.22"
$ws.Range("E8").Value = "  -1.52%  "
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "'0.177"
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("D11").Value = "'0.659"
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("D12").Value = "'55.86"
$ws.Range("E12").Value = "  -3.40%  "
$ws.Range("D13").Value = "'0.0000302"
$ws.Range("E13").Value = "  +2.13%  "
$ws.Range("D14").Value = "'9.57"
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("D15").Value = "4.154.48"
$ws.Range("E15").Value = "  -1.48%  "
$ws.Range("D16").Value = "'19.82"
$ws.Range("E16").Value = "  +2.14%  "
$ws.Range("D17").Value = "3.580.75"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").Value = "69.900.86"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("D19").Value = "'12.56"
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("E21").Value = "  -1.26%  "
$ws.Range("D22").Value = "'474.69"
$ws.Range("E22").Value = "  -3.64%  "
$ws.Range("D23").Value = "'19.07"
$ws.Range("E23").Value = "  +13.57%  "
$ws.Range("D24").Value = "'5.05"
$ws.Range("E24").Value = "  -8.29%  "
$ws.Range("D25").Value = "'4.35"
$ws.Range("E25").Value = "  -2.66%  "
$ws.Range("D26").Value = "'93.30"
$ws.Range("E26").Value = "  +2.83%  "
$ws.Range("D27").Value = "'3.02"
$ws.Range("E27").Value = "  -2.50%  "
$ws.Range("D28").Value = "'10.97"
$ws.Range("E28").Value = "  -2.47%  "
$ws.Range("D29").Value = "'9.29"
$ws.Range("E29").Value = "  -1.11%  "
$ws.Range("D30").Value = "'32.18"
$ws.Range("E30").Value = "  -0.57%  "
$ws.Range("D31").Value = "'7.72"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").Value = "'0.119"
$ws.Range("E32").Value = "  +2.13%  "
$ws.Range("D33").Value = "'12.16"
$ws.Range("E33").Value = "  -0.81%  "
$ws.Range("D34").Value = "'66.16"
$ws.Range("E34").Value = "  +0.88%  "
$ws.Range("D35").Value = "'582.24"
$ws.Range("E35").Value = "  -5.23%  "
$ws.Range("D36").Value = "'38.98"
$ws.Range("E36").Value = "  +2.80%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").Value = "0.0₃0796"
$ws.Range("E38").Value = "  -4.27%  "
$ws.Range("D39").Value = "'0.395"
$ws.Range("E39").Value = "  -2.04%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.139"
$ws.Range("E40").Value = "  -6.38%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'3.22"
$ws.Range("E41").Value = "  +17.78%  "
$ws.Range("D42").Value = "'3.47"
$ws.Range("E42").Value = "  -5.04%  "
$ws.Range("D43").Value = "3.228.12"
$ws.Range("E43").Value = "  -4.00%  "
$ws.Range("D44").Value = "'2.84"
$ws.Range("E44").Value = "  +6.65%  "
$ws.Range("D45").Value = "'3.08"
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("D46").Value = "'0.0442"
$ws.Range("E46").Value = "  -1.05%  "
$ws.Range("D47").Value = "'9.50"
$ws.Range("E47").Value = "  +3.12%  "
$ws.Range("D48").Value = "'3.34"
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").Value = "'3.14"
$ws.Range("E51").Value = "  -6.33%  "

# Reset number-format/style side effects from the quote-prefix (text-forcing) trick
# above so the affected cells keep the workbook-default style, matching the source.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
